# Update Masamune Profits leve tables across sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 662.8
$ws.Range("I92").Value = 717.41174
$ws.Range("J92").Value = 353.33334
$ws.Range("K92").Value = 717.41174
$ws.Range("L92").Value = 353.33334
$ws.Range("M92").Value = 530.58826
$ws.Range("N92").Value = -2849.33334

$ws.Range("H96").Value = 66741020
$ws.Range("I96").Value = 5103.7144
$ws.Range("J96").Value = 125134940
$ws.Range("K96").Value = 15311.1432
$ws.Range("L96").Value = 375404820
$ws.Range("M96").Value = -13938.1432
$ws.Range("N96").Value = -375407566

$ws.Range("H99").Value = 1132
$ws.Range("I99").Value = 1132
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3396
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1898
$ws.Range("N99").ClearContents()

$ws.Range("H101").Value = 535.375
$ws.Range("I101").Value = 306.75
$ws.Range("J101").Value = 1221.25
$ws.Range("K101").Value = 920.25
$ws.Range("L101").Value = 3663.75
$ws.Range("M101").Value = 701.75
$ws.Range("N101").Value = -6907.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37438.75
$ws.Range("I32").Value = 37104.66
$ws.Range("K32").Value = 37104.66
$ws.Range("M32").Value = -36817.66

$ws.Range("H97").Value = 1166.55
$ws.Range("I97").Value = 955.38464
$ws.Range("J97").Value = 1558.7142
$ws.Range("K97").Value = 955.38464
$ws.Range("L97").Value = 1558.7142
$ws.Range("M97").Value = -459.38464
$ws.Range("N97").Value = -2550.7142

$ws.Range("H122").Value = 3011.3845
$ws.Range("I122").Value = 3137.3333
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 9411.999899999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -6961.999899999999
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 916
$ws.Range("I94").Value = 877.5454999999999
$ws.Range("J94").Value = 1000.6
$ws.Range("K94").Value = 877.5454999999999
$ws.Range("L94").Value = 1000.6
$ws.Range("M94").Value = -426.5454999999999
$ws.Range("N94").Value = -1902.6

$ws.Range("H122").Value = 40538.668
$ws.Range("J122").Value = 40538.668
$ws.Range("L122").Value = 40538.668
$ws.Range("N122").Value = -50338.668

$ws.Range("H134").Value = 3755.8674
$ws.Range("I134").Value = 2823.5264
$ws.Range("J134").Value = 4032.6562
$ws.Range("K134").Value = 8470.5792
$ws.Range("L134").Value = 12097.9686
$ws.Range("M134").Value = -5935.5792
$ws.Range("N134").Value = -17167.9686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 36792.89
$ws.Range("I59").Value = 25052
$ws.Range("J59").Value = 40147.43
$ws.Range("K59").Value = 25052
$ws.Range("L59").Value = 40147.43
$ws.Range("M59").Value = -23907
$ws.Range("N59").Value = -42437.43

$ws.Range("H110").Value = 41696.8
$ws.Range("J110").Value = 41696.8
$ws.Range("L110").Value = 41696.8
$ws.Range("N110").Value = -49876.8

$ws.Range("H111").Value = 42549
$ws.Range("J111").Value = 42549
$ws.Range("L111").Value = 42549
$ws.Range("N111").Value = -50729

$ws.Range("H112").Value = 32017.8
$ws.Range("J112").Value = 32017.8
$ws.Range("L112").Value = 32017.8
$ws.Range("N112").Value = -34971.8

$ws.Range("H115").Value = 28922.334
$ws.Range("J115").Value = 28922.334
$ws.Range("L115").Value = 28922.334
$ws.Range("N115").Value = -31272.334

$ws.Range("H116").Value = 40918.75
$ws.Range("J116").Value = 40918.75
$ws.Range("L116").Value = 40918.75
$ws.Range("N116").Value = -50096.75

$ws.Range("H122").Value = 55391.773
$ws.Range("I122").Value = 109690.45
$ws.Range("J122").Value = 1093.091
$ws.Range("K122").Value = 329071.35
$ws.Range("L122").Value = 3279.273
$ws.Range("M122").Value = -326621.35
$ws.Range("N122").Value = -8179.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50000120
$ws.Range("I2").Value = 51.25
$ws.Range("J2").Value = 107143050
$ws.Range("K2").Value = 307.5
$ws.Range("L2").Value = 642858300
$ws.Range("M2").Value = -194.5
$ws.Range("N2").Value = -642858526

$ws.Range("H100").Value = 5294.1113
$ws.Range("J100").Value = 5294.1113
$ws.Range("L100").Value = 15882.3339
$ws.Range("N100").Value = -17504.3339

$ws.Range("H107").Value = 9925.190000000001
$ws.Range("J107").Value = 14813.857
$ws.Range("L107").Value = 44441.571
$ws.Range("N107").Value = -48281.571

$ws.Range("H116").Value = 1243
$ws.Range("I116").Value = 1243
$ws.Range("K116").Value = 3729
$ws.Range("M116").Value = -287

$ws.Range("H121").Value = 282240.8
$ws.Range("I121").Value = 266
$ws.Range("J121").Value = 517219.84
$ws.Range("K121").Value = 798
$ws.Range("L121").Value = 1551659.52
$ws.Range("M121").Value = 512
$ws.Range("N121").Value = -1554279.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1603.8462
$ws.Range("I122").Value = 1468.1818
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 4404.5454
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -1954.5454
$ws.Range("N122").Value = -11950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 46807.332
$ws.Range("J36").Value = 46807.332
$ws.Range("L36").Value = 46807.332
$ws.Range("N36").Value = -47931.332

$ws.Range("H93").Value = 2214
$ws.Range("I93").Value = 1756.8334
$ws.Range("J93").Value = 2442.5833
$ws.Range("K93").Value = 1756.8334
$ws.Range("L93").Value = 2442.5833
$ws.Range("M93").Value = -508.8334
$ws.Range("N93").Value = -4938.5833

$ws.Range("H119").Value = 47408
$ws.Range("J119").Value = 47408
$ws.Range("L119").Value = 47408
$ws.Range("N119").Value = -57084

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45928.25
$ws.Range("J16").Value = 45928.25
$ws.Range("L16").Value = 45928.25
$ws.Range("N16").Value = -46512.25

$ws.Range("H100").Value = 467.6
$ws.Range("I100").Value = 401.16666
$ws.Range("J100").Value = 733.3333
$ws.Range("K100").Value = 802.33332
$ws.Range("L100").Value = 1466.6666
$ws.Range("M100").Value = -261.33332
$ws.Range("N100").Value = -2548.6666

$ws.Range("H122").Value = 2198970.5
$ws.Range("I122").Value = 3175857.2
$ws.Range("J122").Value = 975
$ws.Range("K122").Value = 9527571.600000001
$ws.Range("L122").Value = 2925
$ws.Range("M122").Value = -9525121.600000001
$ws.Range("N122").Value = -7825
